$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("How many tracks can you define in one ODF?", "llama3.2:latest", "According to the Track Settings dialog box, the number of tracks that can be defined is 200."),
    @("How many curve shades can I create?", "llama3.2:latest", "According to the document, you can create 250 curve shades."),
    @("How many curves can I load in one go?", "llama3.2:latest", "You can load up to 450 curves at a time."),
    @("What the maximum number of headers I can display in my log?", "llama3.2:latest", "The maximum number of headers you can display in your log is up to 50."),
    @("How many tables can I have in my log?", "llama3.2:latest", "You can have up to 100 tables in a log."),
    @("How many symbols can I have in the plot at any one time?", "llama3.2:latest", "You can have up to 10,000 symbols in a plot at any given time.")
)

$startRow = 23
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
